$wb = $excel.ActiveWorkbook

# OFF sheet - Week 16 Home row (row 2) target depth data
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 155
$wsOff.Range("C2").Value = 112
$wsOff.Range("D2").Value = 46
$wsOff.Range("E2").Value = 23

# DEF sheet - Week 16 Home row (row 2) target depth data
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 302
$wsDef.Range("C2").Value = 213
$wsDef.Range("D2").Value = 53
$wsDef.Range("E2").Value = 23
